$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 becomes "max" (string), B2/C2 swap, D2 stays "knn"
$ws.Range("A2").Value = "max"
$ws.Range("B2").Value = "svr"
$ws.Range("C2").Value = "knn"
$ws.Range("D2").Value = "knn"

# Row 3: A3 becomes "mean", B3/C3/D3 become strings
$ws.Range("A3").Value = "mean"
$ws.Range("B3").Value = "krr"
$ws.Range("C3").Value = "knn"
$ws.Range("D3").Value = "knn"

# Row 4: new row
$ws.Range("A4").Value = "min"
$ws.Range("B4").Value = "krr"
$ws.Range("C4").Value = "knn"
$ws.Range("D4").Value = "svr"

# Copy style from A3 (style index 1: bold, border, center) to A4
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
